$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure changed cells keep their original Text format (values are
# stored as plain text in the source data, e.g. "43.187.57", "0.500").
$ws.Range('D2').NumberFormat = '@'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('B42').NumberFormat = '@'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('B43').NumberFormat = '@'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('E51').NumberFormat = '@'

$ws.Range('D2').Value = '43.187.57'
$ws.Range('E2').Value = '  +1.80%  '
$ws.Range('D3').Value = '2.378.32'
$ws.Range('E3').Value = '  +3.80%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = '303.25'
$ws.Range('E5').Value = '  +0.95%  '
$ws.Range('D6').Value = '97.33'
$ws.Range('E6').Value = '  +2.66%  '
$ws.Range('D7').Value = '0.508'
$ws.Range('E7').Value = '  +0.40%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').Value = '0.500'
$ws.Range('E9').Value = '  +2.12%  '
$ws.Range('D10').Value = '34.35'
$ws.Range('E10').Value = '  +0.03%  '
$ws.Range('D11').Value = '0.0788'
$ws.Range('E11').Value = '  +1.13%  '
$ws.Range('E12').Value = '  +2.32%  '
$ws.Range('D13').Value = '18.35'
$ws.Range('E13').Value = '  -3.00%  '
$ws.Range('E14').Value = '  +2.03%  '
$ws.Range('D15').Value = '2.752.73'
$ws.Range('E15').Value = '  +4.03%  '
$ws.Range('D16').Value = '2.373.25'
$ws.Range('E16').Value = '  +4.25%  '
$ws.Range('E17').Value = '  +4.14%  '
$ws.Range('D18').Value = '43.185.23'
$ws.Range('E18').Value = '  +1.90%  '
$ws.Range('D19').Value = '12.18'
$ws.Range('E19').Value = '  +0.20%  '
$ws.Range('D20').Value = '6.31'
$ws.Range('E20').Value = '  +5.66%  '
$ws.Range('D21').Value = '0.0₃0889'
$ws.Range('E21').Value = '  +0.39%  '
$ws.Range('D22').Value = '68.48'
$ws.Range('E22').Value = '  +1.50%  '
$ws.Range('D23').Value = '235.64'
$ws.Range('E23').Value = '  +0.22%  '
$ws.Range('D24').Value = '2.23'
$ws.Range('E24').Value = '  -0.44%  '
$ws.Range('E25').Value = '  +2.44%  '
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('D27').Value = '24.83'
$ws.Range('E27').Value = '  +2.85%  '
$ws.Range('E28').Value = '  +0.40%  '
$ws.Range('D29').Value = '9.15'
$ws.Range('E29').Value = '  +1.56%  '
$ws.Range('D30').Value = '31.71'
$ws.Range('E30').Value = '  +0.49%  '
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('E32').Value = '  +2.75%  '
$ws.Range('D33').Value = '0.0740'
$ws.Range('E33').Value = '  +7.17%  '
$ws.Range('D34').Value = '17.29'
$ws.Range('E34').Value = '  -0.54%  '
$ws.Range('E35').Value = '  +5.47%  '
$ws.Range('E36').Value = '  +7.20%  '
$ws.Range('E37').Value = '  -1.07%  '
$ws.Range('E38').Value = '  -1.55%  '
$ws.Range('E39').Value = '  +4.57%  '
$ws.Range('D40').Value = '22.55'
$ws.Range('E40').Value = '  +14.59%  '
$ws.Range('D41').Value = '0.108'
$ws.Range('E41').Value = '  +0.44%  '
$ws.Range('B42').Value = 'Monero'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D42').Value = '104.58'
$ws.Range('E42').Value = '  -36.40%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '1.958.92'
$ws.Range('E43').Value = '  +0.78%  '
$ws.Range('E44').Value = '  +1.34%  '
$ws.Range('E45').Value = '  +2.36%  '
$ws.Range('D46').Value = '2.75'
$ws.Range('E46').Value = '  +1.22%  '
$ws.Range('D47').Value = '9.20'
$ws.Range('E47').Value = '  -10.36%  '
$ws.Range('D48').Value = '52.77'
$ws.Range('E48').Value = '  +0.19%  '
$ws.Range('E49').Value = '  +3.39%  '
$ws.Range('D50').Value = '71.92'
$ws.Range('E50').Value = '  +2.43%  '
$ws.Range('D51').Value = '1.14'
$ws.Range('E51').Value = '  +1.84%  '
